# fix(publipostage): Refactor synthetic array /3
# Rename the shared-string "color" tokens used in the intervention_type
# legend columns (A = emoji swatch, B = color name):
#   black square ⬛  -> blue book  📘
#   red square   🟥  -> red book   📕
#   orange square🟧  -> orange book📙
#   green square 🟩  -> green book 📗
#   noir             -> bleu

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlWhole = 1, xlPart = 2 ; xlByRows = 1, xlByColumns = 2
$xlWhole = 1
$xlByRows = 1

$cells = $ws.Range("A1:M121")

$cells.Replace("⬛", "📘", $xlWhole, $xlByRows, $false, $false, $true)
$cells.Replace("🟥", "📕", $xlWhole, $xlByRows, $false, $false, $true)
$cells.Replace("🟧", "📙", $xlWhole, $xlByRows, $false, $false, $true)
$cells.Replace("🟩", "📗", $xlWhole, $xlByRows, $false, $false, $true)
$cells.Replace("noir", "bleu", $xlWhole, $xlByRows, $false, $false, $true)
